# Re-shuffle the per-trial stimulus records (category/condition/correct-answer
# and the stimulus file + conceptual/perceptual/typicality stats columns)
# across rows 2-41, leaving the trial-structure columns (A-G, J) untouched.
#
# This reproduces the commit's effect of deriving a new "subject" ordering
# from the same pool of stimuli (one of only 20 distinct underlying
# orderings that get duplicated across the 1000 generated subject files).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel together as one "stimulus record".
$cols = @(8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22)

# Destination row -> source row (which row's pre-edit record lands here).
$rowMap = @{2=20; 3=16; 4=12; 5=40; 6=39; 7=33; 8=2; 9=23; 10=15; 11=14; 12=35; 13=5; 14=21; 15=29; 16=24; 17=36; 18=30; 19=18; 20=10; 21=13; 22=37; 23=11; 24=3; 25=41; 26=38; 27=4; 28=22; 29=25; 30=17; 31=31; 32=32; 33=8; 34=26; 35=19; 36=28; 37=7; 38=27; 39=6; 40=9; 41=34}

# 1) Snapshot every row's current ("before") record.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rec = @{}
    foreach ($c in $cols) {
        $rec[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rec
}

# 2) Write each row's new record using the source row's snapshotted values.
for ($r = 2; $r -le 41; $r++) {
    $srcRow = $rowMap[$r]
    $srcRec = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcRec[$c]
    }
}
